# Generate Report for Handback
# Update the timestamps / status recorded on the handback report:
#  - "Latest HO Xliff Generate Date" (Overview) / "Correspond Handoff Datetime" (de-de)
#  - "Priority" column (zh-cn / de-de) ht -> mt
#  - "Correspond Handoff Datetime" / "Correspond Handback DateTime" (zh-cn)
#  - "Correspond Handback DateTime" (de-de)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column G, rows 2 and 5
$wsOverview.Range("G2").Value = "2016-08-24 14:17:33"
$wsOverview.Range("G5").Value = "2016-08-24 14:17:33"

# zh-cn sheet: Priority column E, rows 2 and 5 (ht -> mt)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime column H, rows 2 and 5
$wsZhCn.Range("H2").Value = "2016-08-24 14:17:27"
$wsZhCn.Range("H5").Value = "2016-08-24 14:17:27"

# zh-cn sheet: Correspond Handback DateTime column K, rows 2 and 5
$wsZhCn.Range("K2").Value = "2016-08-24 14:17:44"
$wsZhCn.Range("K5").Value = "2016-08-24 14:17:44"

# de-de sheet: Priority column E, rows 2 and 5 (ht -> mt)
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# de-de sheet: Correspond Handoff Datetime column H, rows 2 and 5
$wsDeDe.Range("H2").Value = "2016-08-24 14:17:33"
$wsDeDe.Range("H5").Value = "2016-08-24 14:17:33"

# de-de sheet: Correspond Handback DateTime column K, rows 2 and 5
$wsDeDe.Range("K2").Value = "2016-08-24 14:17:51"
$wsDeDe.Range("K5").Value = "2016-08-24 14:17:51"
